# Apply the corrected "Fitness" values (column C) for run_9.xlsx
# as described by the commit "correction in sa algorithm and 746 logs".
# The new values collapse into a small number of contiguous runs,
# so we can set each contiguous block of rows in one shot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    @{ Start = 2;   End = 31;  Value = 7859 },
    @{ Start = 32;  End = 72;  Value = 7812 },
    @{ Start = 73;  End = 80;  Value = 7785 },
    @{ Start = 81;  End = 84;  Value = 7320 },
    @{ Start = 85;  End = 87;  Value = 7318 },
    @{ Start = 88;  End = 252; Value = 7310 }
)

foreach ($r in $ranges) {
    $rangeAddress = "C" + $r.Start + ":C" + $r.End
    $ws.Range($rangeAddress).Value = $r.Value
}
